$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
$tcs = $np.ThemeColorScheme
Write-Host "count:" $tcs.Count
for ($i=1; $i -le $tcs.Count; $i++) {
  Write-Host $i $tcs.Colors($i).RGB
}
